$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: updated electricity spot price data (automated daily price update)
$ws.Range("A2").Value = 45914
$ws.Range("B2").Value = 75
$ws.Range("C2").Value = 66.23
$ws.Range("D2").Value = 65.29000000000001
$ws.Range("E2").Value = 65.3
$ws.Range("F2").Value = 65.81
$ws.Range("G2").Value = 71.54000000000001
$ws.Range("H2").Value = 71.54000000000001
$ws.Range("I2").Value = 77.56999999999999
$ws.Range("J2").Value = 66.23
$ws.Range("K2").Value = 23.58
$ws.Range("L2").Value = 5.01
$ws.Range("M2").Value = 0
$ws.Range("P2").Value = -0.01
$ws.Range("S2").Value = 1.1
$ws.Range("T2").Value = 15.65
$ws.Range("U2").Value = 67.3
$ws.Range("V2").Value = 103.81
$ws.Range("W2").Value = 126.13
$ws.Range("X2").Value = 106.01
$ws.Range("Y2").Value = 99.84
$ws.Range("Z2").Value = 48.87
$ws.Range("AB2").Value = 108.95
$ws.Range("AD2").Value = 114.97
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 102.93
$ws.Range("AG2").Value = "9h-18h"
